$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1662531017369727
$ws.Range("C2").Value = 0.5880893300248139
$ws.Range("J2").Value = 0.01488833746898263
$ws.Range("P2").Value = 0.1240694789081886
$ws.Range("S2").Value = 0.1066997518610422
$ws.Range("B3").Value = 0.00816326530612245
$ws.Range("C3").Value = 0.0163265306122449
$ws.Range("J3").Value = 0.04081632653061224
$ws.Range("P3").Value = 0.7346938775510204
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.01754385964912281
$ws.Range("P4").Value = 0.7368421052631579
$ws.Range("S4").Value = 0.2456140350877193
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.04721030042918455
$ws.Range("D6").Value = 0.01716738197424893
$ws.Range("E6").Value = 0.008583690987124463
$ws.Range("F6").Value = 0.03004291845493562
$ws.Range("J6").Value = 0.3218884120171674
$ws.Range("O6").Value = 0.01716738197424893
$ws.Range("Q6").Value = 0.1802575107296137
$ws.Range("R6").Value = 0.04291845493562232
$ws.Range("S6").Value = 0.3347639484978541
$ws.Range("B7").Value = 0.1191709844559585
$ws.Range("D7").Value = 0.0155440414507772
$ws.Range("F7").Value = 0.02590673575129534
$ws.Range("J7").Value = 0.2020725388601036
$ws.Range("O7").Value = 0.0155440414507772
$ws.Range("Q7").Value = 0.2227979274611399
$ws.Range("R7").Value = 0.0310880829015544
$ws.Range("S7").Value = 0.3678756476683938
$ws.Range("B8").Value = 0.1383219954648526
$ws.Range("D8").Value = 0.02267573696145125
$ws.Range("E8").Value = 0.00453514739229025
$ws.Range("F8").Value = 0.04988662131519275
$ws.Range("J8").Value = 0.1020408163265306
$ws.Range("O8").Value = 0.02040816326530612
$ws.Range("Q8").Value = 0.1972789115646258
$ws.Range("R8").Value = 0.08843537414965986
$ws.Range("S8").Value = 0.3764172335600907
$ws.Range("B9").Value = 0.1325966850828729
$ws.Range("D9").Value = 0.005524861878453038
$ws.Range("E9").Value = 0.005524861878453038
$ws.Range("F9").Value = 0.03867403314917127
$ws.Range("J9").Value = 0.1270718232044199
$ws.Range("O9").Value = 0.02209944751381215
$ws.Range("Q9").Value = 0.1988950276243094
$ws.Range("R9").Value = 0.08287292817679558
$ws.Range("S9").Value = 0.3867403314917127
$ws.Range("B10").Value = 0.1341935483870968
$ws.Range("D10").Value = 0.02580645161290323
$ws.Range("F10").Value = 0.06903225806451613
$ws.Range("J10").Value = 0.1316129032258065
$ws.Range("O10").Value = 0.02
$ws.Range("Q10").Value = 0.2161290322580645
$ws.Range("R10").Value = 0.06774193548387097
$ws.Range("S10").Value = 0.335483870967742
$ws.Range("G11").Value = 0.125
$ws.Range("J11").Value = 0.1121794871794872
$ws.Range("K11").Value = 0.1826923076923077
$ws.Range("L11").Value = 0.5673076923076923
$ws.Range("S11").Value = 0.01282051282051282
$ws.Range("G12").Value = 0.7277777777777777
$ws.Range("J12").Value = 0.2333333333333333
$ws.Range("K12").Value = 0.005555555555555556
$ws.Range("L12").Value = 0.01111111111111111
$ws.Range("S12").Value = 0.02222222222222222
$ws.Range("G13").Value = 0.5869565217391305
$ws.Range("J13").Value = 0.391304347826087
$ws.Range("S13").Value = 0.02173913043478261
$ws.Range("F15").Value = 0.02068965517241379
$ws.Range("H15").Value = 0.1241379310344828
$ws.Range("I15").Value = 0.07931034482758621
$ws.Range("J15").Value = 0.3931034482758621
$ws.Range("K15").Value = 0.0896551724137931
$ws.Range("M15").Value = 0.006896551724137931
$ws.Range("O15").Value = 0.06896551724137931
$ws.Range("S15").Value = 0.2172413793103448
$ws.Range("F16").Value = 0.01153846153846154
$ws.Range("H16").Value = 0.1576923076923077
$ws.Range("I16").Value = 0.1076923076923077
$ws.Range("J16").Value = 0.3807692307692307
$ws.Range("K16").Value = 0.1038461538461539
$ws.Range("M16").Value = 0.02692307692307692
$ws.Range("O16").Value = 0.08076923076923077
$ws.Range("S16").Value = 0.1307692307692308
$ws.Range("F17").Value = 0.02587800369685767
$ws.Range("H17").Value = 0.1866913123844732
$ws.Range("I17").Value = 0.05730129390018484
$ws.Range("J17").Value = 0.4547134935304991
$ws.Range("K17").Value = 0.08133086876155268
$ws.Range("M17").Value = 0.022181146025878
$ws.Range("N17").Value = 0.001848428835489834
$ws.Range("O17").Value = 0.07024029574861368
$ws.Range("S17").Value = 0.09981515711645102
$ws.Range("F18").Value = 0.01734104046242774
$ws.Range("H18").Value = 0.1560693641618497
$ws.Range("I18").Value = 0.1040462427745665
$ws.Range("J18").Value = 0.4624277456647399
$ws.Range("K18").Value = 0.07514450867052024
$ws.Range("M18").Value = 0.005780346820809248
$ws.Range("O18").Value = 0.08670520231213873
$ws.Range("S18").Value = 0.09248554913294797
$ws.Range("F19").Value = 0.02515243902439025
$ws.Range("H19").Value = 0.1814024390243902
$ws.Range("I19").Value = 0.06326219512195122
$ws.Range("J19").Value = 0.4047256097560976
$ws.Range("K19").Value = 0.1089939024390244
$ws.Range("M19").Value = 0.02057926829268293
$ws.Range("O19").Value = 0.08460365853658537
$ws.Range("S19").Value = 0.1112804878048781
